$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data values (B:G) for rows 2-9, replacing the old s_vals dataset
# after filtering out save games.
$data = @{
    2 = @(1.505614041169197, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054)
    3 = @(3.182878228561681, 1.65323645889881,  3.082599426703578,  0.4998867070740569, 1, 8.418600821238126)
    4 = @(3.182878228561681, 1.65323645889881,  3.082599426703578,  0.4998867070740569, 0, 8.418600821238126)
    5 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    6 = @(0.7287194209349384,1.65323645889881,  0.7127328510149897, 0.4998867070740569, 0, 3.594575437922795)
    7 = @(0.3464964993005633,0.3375848360084654,0.1529057820181812, 0.4998867070740569, 0, 1.336873824401267)
    8 = @(3.182878228561681, 1.65323645889881,  16.98373111632243,  0.4998867070740569, 0, 22.31973251085698)
    9 = @(0.7287194209349384,1.65323645889881,  0.7127328510149897, 0.4998867070740569, 0, 3.594575437922795)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
